{"js": "// Update the header date and the 25 two-digit-by-two-digit multiplication\n// prompts to the new values from the commit.\nconst replacements = [\n  [\"2025-12-29 Monday\", \"2025-12-30 Tuesday\"],\n  [\"23\u00d713=\", \"48\u00d721=\"],\n  [\"92\u00d792=\", \"74\u00d796=\"],\n  [\"49\u00d715=\", \"58\u00d787=\"],\n  [\"36\u00d730=\", \"80\u00d787=\"],\n  [\"46\u00d747=\", \"13\u00d765=\"],\n  [\"81\u00d742=\", \"33\u00d733=\"],\n  [\"91\u00d786=\", \"11\u00d762=\"],\n  [\"39\u00d748=\", \"35\u00d720=\"],\n  [\"16\u00d715=\", \"62\u00d770=\"],\n  [\"92\u00d799=\", \"46\u00d778=\"],\n  [\"32\u00d782=\", \"84\u00d785=\"],\n  [\"54\u00d750=\", \"69\u00d757=\"],\n  [\"37\u00d780=\", \"90\u00d749=\"],\n  [\"67\u00d720=\", \"91\u00d769=\"],\n  [\"58\u00d788=\", \"69\u00d773=\"],\n  [\"24\u00d786=\", \"79\u00d722=\"],\n  [\"14\u00d767=\", \"56\u00d732=\"],\n  [\"80\u00d712=\", \"26\u00d762=\"],\n  [\"91\u00d731=\", \"52\u00d721=\"],\n  [\"65\u00d782=\", \"49\u00d784=\"],\n  [\"32\u00d753=\", \"51\u00d737=\"],\n  [\"20\u00d733=\", \"75\u00d781=\"],\n  [\"27\u00d788=\", \"97\u00d722=\"],\n  [\"72\u00d753=\", \"20\u00d741=\"],\n  [\"65\u00d790=\", \"88\u00d760=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the header date and the 25 two-digit-by-two-digit multiplication\n# prompts to the new values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-12-29 Monday\", \"2025-12-30 Tuesday\"),\n    @(\"23\u00d713=\", \"48\u00d721=\"),\n    @(\"92\u00d792=\", \"74\u00d796=\"),\n    @(\"49\u00d715=\", \"58\u00d787=\"),\n    @(\"36\u00d730=\", \"80\u00d787=\"),\n    @(\"46\u00d747=\", \"13\u00d765=\"),\n    @(\"81\u00d742=\", \"33\u00d733=\"),\n    @(\"91\u00d786=\", \"11\u00d762=\"),\n    @(\"39\u00d748=\", \"35\u00d720=\"),\n    @(\"16\u00d715=\", \"62\u00d770=\"),\n    @(\"92\u00d799=\", \"46\u00d778=\"),\n    @(\"32\u00d782=\", \"84\u00d785=\"),\n    @(\"54\u00d750=\", \"69\u00d757=\"),\n    @(\"37\u00d780=\", \"90\u00d749=\"),\n    @(\"67\u00d720=\", \"91\u00d769=\"),\n    @(\"58\u00d788=\", \"69\u00d773=\"),\n    @(\"24\u00d786=\", \"79\u00d722=\"),\n    @(\"14\u00d767=\", \"56\u00d732=\"),\n    @(\"80\u00d712=\", \"26\u00d762=\"),\n    @(\"91\u00d731=\", \"52\u00d721=\"),\n    @(\"65\u00d782=\", \"49\u00d784=\"),\n    @(\"32\u00d753=\", \"51\u00d737=\"),\n    @(\"20\u00d733=\", \"75\u00d781=\"),\n    @(\"27\u00d788=\", \"97\u00d722=\"),\n    @(\"72\u00d753=\", \"20\u00d741=\"),\n    @(\"65\u00d790=\", \"88\u00d760=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $oldText\n    $rng.Find.Replacement.Text = $newText\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\n$d.Save()\n"}
